$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary field updates ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 35

# --- Step 1: relocate the TOTAL row formatting + merge from old row 39 to new row 51 ---
# (must run before the zebra-stripe extension below re-touches row 39)
$ws.Range("A39:H39").Copy()
$ws.Range("A51:H51").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A39:G39").UnMerge()
$ws.Range("A51:G51").Merge()

# --- Step 2: extend the alternating (zebra-stripe) row shading down through row 50 ---
# Row 16 carries the "even" style set (9,10,11); row 17 carries the "odd" style set (12,13,14).
# PasteSpecial only honours the first area of a multi-area range in this host, so loop per row.
for ($r = 18; $r -le 50; $r++) {
    if ($r % 2 -eq 0) {
        $srcRow = 16
    } else {
        $srcRow = 17
    }
    $ws.Range("A" + $srcRow + ":H" + $srcRow).Copy()
    $ws.Range("A" + $r + ":H" + $r).PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false

# --- Step 3: write the detail-table content (rows 16-50) ---
$rows = @(
    @{A="Point 09"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 17"; B="GND-CR-4"; C="Inst"; D="GND,Cu Clad Rod,#4"; E="EA"; F=1; H=0.0}
    @{A="Point 17"; B="GND-MD"; C="Inst"; D="GND,Wire Mldg Only"; E="EA"; F=2; H=0.0}
    @{A="Point 17"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 19"; B="GND-LG-SL"; C="Inst"; D="GND,Lug,St Lt"; E="EA"; F=1; H=0.0}
    @{A="Point 19"; B="GND-MD"; C="Inst"; D="GND,Wire Mldg Only"; E="EA"; F=2; H=0.0}
    @{A="Point 19"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 21"; B="INS-15-P-S-C"; C="Inst"; D="INS,15kV,Pin,Silicon Polymer,Corr"; E="EA"; F=1; H=0.0}
    @{A="Point 21"; B="PIN-15-PTP-C"; C="Inst"; D="Pin,15kV,Pole top,Corrosive"; E="EA"; F=1; H=0.0}
    @{A="Point 21"; B="POL-40-2"; C="Inst"; D="Pole,40ft,Class 2"; E="EA"; F=1; H=0.0}
    @{A="Point 21"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 23"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="GND-CR-4"; C="Inst"; D="GND,Cu Clad Rod,#4"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="GND-LG-SL"; C="Inst"; D="GND,Lug,St Lt"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="GND-MD"; C="Inst"; D="GND,Wire Mldg Only"; E="EA"; F=2; H=0.0}
    @{A="Point 25"; B="INS-15-D-S-C"; C="Inst"; D="INS,15kV,Deadend,Polymer,Corr"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="INS-15-P-S-C"; C="Inst"; D="INS,15kV,Pin,Silicon Polymer,Corr"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="PIN-15-PTP-C"; C="Inst"; D="Pin,15kV,Pole top,Corrosive"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="POL-40-2"; C="Inst"; D="Pole,40ft,Class 2"; E="EA"; F=1; H=0.0}
    @{A="Point 25"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 27"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 31"; B="GND-CR-4"; C="Inst"; D="GND,Cu Clad Rod,#4"; E="EA"; F=1; H=0.0}
    @{A="Point 31"; B="GND-MD"; C="Inst"; D="GND,Wire Mldg Only"; E="EA"; F=2; H=0.0}
    @{A="Point 31"; B="INS-15-D-S-C"; C="Inst"; D="INS,15kV,Deadend,Polymer,Corr"; E="EA"; F=1; H=0.0}
    @{A="Point 31"; B="PLD-EYE-C"; C="Inst"; D="PLD,Eyebolt Deadend,Corrosive"; E="EA"; F=4; H=0.0}
    @{A="Point 31"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 33"; B="INS-15-D-S-C"; C="Inst"; D="INS,15kV,Deadend,Polymer,Corr"; E="EA"; F=1; H=0.0}
    @{A="Point 33"; B="PLD-EYE-C"; C="Inst"; D="PLD,Eyebolt Deadend,Corrosive"; E="EA"; F=2; H=0.0}
    @{A="Point 33"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
    @{A="Point 37"; B="GND-CR-4"; C="Inst"; D="GND,Cu Clad Rod,#4"; E="EA"; F=1; H=0.0}
    @{A="Point 37"; B="GND-MD"; C="Inst"; D="GND,Wire Mldg Only"; E="EA"; F=2; H=0.0}
    @{A="Point 37"; B="INS-15-P-S-C"; C="Inst"; D="INS,15kV,Pin,Silicon Polymer,Corr"; E="EA"; F=1; H=0.0}
    @{A="Point 37"; B="PIN-15-PTP-C"; C="Inst"; D="Pin,15kV,Pole top,Corrosive"; E="EA"; F=1; H=0.0}
    @{A="Point 37"; B="POL-40-2"; C="Inst"; D="Pole,40ft,Class 2"; E="EA"; F=1; H=0.0}
    @{A="Point 37"; B="PLA-HDIG"; C="Inst"; D="PLA,Hand Dig or Additional  Excavation"; E="EA"; F=1; H=0.0}
)

$r = 16
foreach ($row in $rows) {
    $ws.Range("A" + $r).Value = $row.A
    $ws.Range("B" + $r).Value = $row.B
    $ws.Range("C" + $r).Value = $row.C
    $ws.Range("D" + $r).Value = $row.D
    $ws.Range("E" + $r).Value = $row.E
    $ws.Range("F" + $r).Value = $row.F
    $ws.Range("H" + $r).Value = $row.H
    $r++
}

# --- Step 4: TOTAL row (51) ---
$ws.Range("A51").Value = "TOTAL"
$ws.Range("H51").Value = 0
